# fix: unificacion de clases
# Renumbers the token codes in the "Tabla de simbolos" sheet and adds a new
# literal/lexema row ("{ Hola mundo }" -> "-> cadena") to the second table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B (Token) numeric updates for the literal/reserved-word table (rows 4-45) ---
$bUpdates = @{
  4=95;   5=43;   6=45;   7=42;   8=47;   9=58;   10=62;  11=60;  12=61;
  13=33;  14=40;  15=41;  16=44;  17=46;  18=59;  19=35;  20=92;  21=123;
  22=125; 23=10;  24=32;  25=9;   26=257; 27=258; 28=259; 29=260; 30=261;
  31=262; 32=263; 33=264; 34=265; 35=266; 36=267; 37=268; 38=269; 39=270;
  40=271; 41=272; 42=273; 43=274; 44=275; 45=276
}
foreach ($row in $bUpdates.Keys) {
  $ws.Cells.Item($row, 2).Value = $bUpdates[$row]
}

# --- Column E (Token) numeric updates for the literal/token/lexema table ---
$eUpdates = @{ 4=280; 7=281; 8=282; 9=282; 10=283; 11=283; 12=281 }
foreach ($row in $eUpdates.Keys) {
  $ws.Cells.Item($row, 5).Value = $eUpdates[$row]
}

# --- New row 13 entry in the literal/token/lexema table: "{ Hola mundo }" / 284 / "-> cadena" ---
# Copy formats from neighbouring rows that already carry the right style (s=4 / s=4 / s=1)
$ws.Range("D12").Copy() | Out-Null
$ws.Range("D13").PasteSpecial(-4122) | Out-Null
$ws.Range("E12").Copy() | Out-Null
$ws.Range("E13").PasteSpecial(-4122) | Out-Null
$ws.Range("F7").Copy() | Out-Null
$ws.Range("F13").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("D13").Value = "{ Hola mundo }"
$ws.Range("E13").Value = 284
$ws.Range("F13").Value = "-> cadena"

# --- Column widths: split B:D into B:C (default 11.52) and a wider D (~13.3) ---
$ws.Range("D1").EntireColumn.ColumnWidth = 12.43

# --- View: scroll back to the top and move the active selection ---
$ws.Range("G17").Select() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 1

Write-Host "edit applied"
